# Update ticket-interest / price numbers across the four sheets of the
# "杭州-漫展信息" workbook (regenerated gh-pages data snapshot).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G2").Value  = 88
$ws.Range("F3").Value  = 627
$ws.Range("F4").Value  = 40
$ws.Range("F5").Value  = 1965
$ws.Range("F6").Value  = 5590
$ws.Range("F7").Value  = 1567
$ws.Range("F9").Value  = 3181
$ws.Range("F12").Value = 1321
$ws.Range("F13").Value = 4436
$ws.Range("F15").Value = 1691
$ws.Range("F18").Value = 41
$ws.Range("F20").Value = 160
$ws.Range("F24").Value = 80
$ws.Range("F25").Value = 10
$ws.Range("F26").Value = 89
$ws.Range("F28").Value = 2
$ws.Range("F31").Value = 69
$ws.Range("F32").Value = 193
$ws.Range("F33").Value = 335
$ws.Range("F35").Value = 11
$ws.Range("F36").Value = 1717
$ws.Range("F37").Value = 2222
$ws.Range("F38").Value = 1038
$ws.Range("F40").Value = 263
$ws.Range("F41").Value = 620
$ws.Range("F42").Value = 334
$ws.Range("F43").Value = 21
$ws.Range("F44").Value = 664
$ws.Range("F45").Value = 18
$ws.Range("F46").Value = 418
$ws.Range("F47").Value = 366
$ws.Range("F48").Value = 221

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 155

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 767

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 767
$ws.Range("G3").Value  = 88
$ws.Range("F4").Value  = 627
$ws.Range("F5").Value  = 40
$ws.Range("F6").Value  = 1965
$ws.Range("F7").Value  = 5590
$ws.Range("F8").Value  = 1567
$ws.Range("F10").Value = 3181
$ws.Range("F12").Value = 1321
$ws.Range("F13").Value = 4436
$ws.Range("F15").Value = 1691
$ws.Range("F17").Value = 41
$ws.Range("F21").Value = 160
$ws.Range("F23").Value = 155
$ws.Range("F26").Value = 80
$ws.Range("F27").Value = 89
$ws.Range("F32").Value = 69
$ws.Range("F33").Value = 193
$ws.Range("F35").Value = 11
$ws.Range("F36").Value = 1717
$ws.Range("F37").Value = 2222
$ws.Range("F38").Value = 1038
$ws.Range("F41").Value = 263
$ws.Range("F42").Value = 620
$ws.Range("F43").Value = 334
$ws.Range("F44").Value = 664
$ws.Range("F45").Value = 418
$ws.Range("F46").Value = 366
$ws.Range("F47").Value = 221
